$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "67.165.01"
$ws.Range("E2").Value = "  +0.94%  "

Set-TextValue "D3" "3.832.81"
$ws.Range("E3").Value = "  +0.89%  "

Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.09%  "

Set-TextValue "D5" "448.21"
$ws.Range("E5").Value = "  +6.84%  "

Set-TextValue "D6" "148.22"
$ws.Range("E6").Value = "  +15.60%  "

Set-TextValue "D7" "0.624"
$ws.Range("E7").Value = "  +4.39%  "

Set-TextValue "D8" "0.999"
$ws.Range("E8").Value = "  -0.05%  "

Set-TextValue "D9" "0.742"
$ws.Range("E9").Value = "  +3.53%  "

$ws.Range("E10").Value = "  -3.10%  "

Set-TextValue "D11" "0.0000324"
$ws.Range("E11").Value = "  -5.37%  "

Set-TextValue "D12" "43.84"
$ws.Range("E12").Value = "  +10.44%  "

Set-TextValue "D13" "10.39"
$ws.Range("E13").Value = "  +4.28%  "

Set-TextValue "D14" "4.439.12"
$ws.Range("E14").Value = "  +1.00%  "

Set-TextValue "D15" "15.03"
$ws.Range("E15").Value = "  -5.67%  "

Set-TextValue "D16" "3.826.41"
$ws.Range("E16").Value = "  +1.34%  "

$ws.Range("E17").Value = "  -0.35%  "

Set-TextValue "D18" "19.98"
$ws.Range("E18").Value = "  +3.31%  "

$ws.Range("E19").Value = "  +7.69%  "

Set-TextValue "D20" "67.263.06"

Set-TextValue "D21" "424.83"
$ws.Range("E21").Value = "  +5.49%  "

$ws.Range("E23").Value = "  +8.79%  "

Set-TextValue "D24" "86.56"
$ws.Range("E24").Value = "  +4.14%  "

Set-TextValue "D25" "37.41"
$ws.Range("E25").Value = "  +1.77%  "

Set-TextValue "D26" "3.44"
$ws.Range("E26").Value = "  +8.67%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D27" "9.60"
$ws.Range("E27").Value = "  +18.60%  "

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D28" "5.49"
$ws.Range("E28").Value = "  -3.81%  "

Set-TextValue "D29" "9.71"
$ws.Range("E29").Value = "  +4.60%  "

Set-TextValue "D30" "750.28"
$ws.Range("E30").Value = "  +6.90%  "

Set-TextValue "D31" "13.74"
$ws.Range("E31").Value = "  +12.52%  "

$ws.Range("E32").Value = "  +12.38%  "

Set-TextValue "D33" "2.74"
$ws.Range("E33").Value = "  -1.22%  "

Set-TextValue "D34" "43.03"
$ws.Range("E34").Value = "  +13.47%  "

$ws.Range("E35").Value = "  +3.79%  "

$ws.Range("E36").Value = "  +6.14%  "

$ws.Range("E37").Value = "  -0.13%  "

Set-TextValue "D38" "5.54"
$ws.Range("E38").Value = "  +17.05%  "

$ws.Range("E39").Value = "  +6.07%  "

Set-TextValue "D40" "0.344"
$ws.Range("E40").Value = "  +18.15%  "

$ws.Range("E41").Value = "  -10.78%  "

$ws.Range("E42").Value = "  -2.77%  "

Set-TextValue "D43" "1.00"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("E44").Value = "  +4.97%  "

Set-TextValue "D45" "3.46"
$ws.Range("E45").Value = "  +4.71%  "

Set-TextValue "D46" "3.24"
$ws.Range("E46").Value = "  +5.67%  "

Set-TextValue "D47" "2.47"
$ws.Range("E47").Value = "  +13.69%  "

$ws.Range("E48").Value = "  +5.08%  "

Set-TextValue "D49" "146.60"
$ws.Range("E49").Value = "  +1.37%  "

Set-TextValue "D50" "2.66"
$ws.Range("E50").Value = "  +5.36%  "

Set-TextValue "D51" "2.88"
$ws.Range("E51").Value = "  +4.97%  "
